$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B, shifting NCTId..results from B..I to C..J
$ws.Columns("B").Insert()

# New header for the inserted column
$ws.Range("B1").Value = "status_label"

# The two clinical-trial rows were also reordered (NCT04560153/rouge now row 2,
# NCT03335345/vert now row 3); update only the cells whose value actually
# changes, leaving the untouched columns (eudraCT, completion_year, ...) as
# the column-insert already shifted them.

# Row 2 becomes the NCT04560153 (red / rouge) trial
$ws.Range("A2").Value = "🟥"
$ws.Range("B2").Value = "rouge"
$ws.Range("C2").Value = "NCT04560153"
$ws.Range("F2").Value = "Interest of Karate Kata Practice on the Self-esteem of Patients Living With HIV"
$ws.Range("G2").Value = "KATACHRO"
$ws.Range("H2").Value = $false
$ws.Range("I2").Value = $false
$ws.Range("J2").Value = $false

# Row 3 becomes the NCT03335345 (green / vert) trial
$ws.Range("A3").Value = "🟩"
$ws.Range("B3").Value = "vert"
$ws.Range("C3").Value = "NCT03335345"
$ws.Range("F3").Value = "Non-inferiority Study of the Pursuit of Enteral Nutrition Compared to a Strategy of Gastric Emptiness Peri-extubation. Cluster Randomized Trial"
$ws.Range("G3").Value = "AMBROISIE"
$ws.Range("H3").Value = $true
$ws.Range("I3").Value = $true
$ws.Range("J3").Value = $true

# Row 4 (NCT03071601) keeps its data, just gains the new status_label value
$ws.Range("B4").Value = "vert"
